# Delete the "Charcoal" row (row 21) from Sheet1, causing all rows below
# it to shift up by one. This matches the diff: the row previously holding
# "Charcoal" (A21) is removed, and the data that used to start at row 22
# ("Convolvulaceae") now starts at row 21, all the way down so what used
# to be row 69 ("Fabaceae/Leguminosa") is now row 68.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sanity check: row 21 column A should currently contain "Charcoal" before
# we delete it (defensive, but harmless if it doesn't match).
$targetRow = 21
$cellValue = $ws.Cells.Item($targetRow, 1).Value

if ($cellValue -eq "Charcoal") {
    $ws.Rows.Item($targetRow).Delete()
} else {
    # Fallback: locate the row containing "Charcoal" in column A and delete it.
    $found = $ws.Range("A1:A100").Find("Charcoal")
    if ($found -ne $null) {
        $found.EntireRow.Delete()
    }
}
